$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New headers for columns AF, AG, AH (row 1) ---
$ws.Range("AF1").Value = "dist_trav_20min_body_out"
$ws.Range("AG1").Value = "dist_trav_25min_body_out"
$ws.Range("AH1").Value = "dist_trav_30min_body_out"

# Match header style (bold + centered) used by the rest of row 1
$headerRange = $ws.Range("AF1:AH1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108

# --- New data values for columns AF, AG, AH (rows 2-78) ---
$ws.Range("AF2").Value = 0
$ws.Range("AG2").Value = 0
$ws.Range("AH2").Value = 0
$ws.Range("AF3").Value = 0
$ws.Range("AG3").Value = 0
$ws.Range("AH3").Value = 0
$ws.Range("AF4").Value = 29.065045753
$ws.Range("AG4").Value = 43.707970517
$ws.Range("AH4").Value = 53.959036928
$ws.Range("AF5").Value = 51.8091856353
$ws.Range("AG5").Value = 0
$ws.Range("AH5").Value = 0
$ws.Range("AF6").Value = 0
$ws.Range("AG6").Value = 0
$ws.Range("AH6").Value = 0
$ws.Range("AF7").Value = 0
$ws.Range("AG7").Value = 0
$ws.Range("AH7").Value = 0
$ws.Range("AF8").Value = 0
$ws.Range("AG8").Value = 0
$ws.Range("AH8").Value = 0
$ws.Range("AF9").Value = 0
$ws.Range("AG9").Value = 0
$ws.Range("AH9").Value = 0
$ws.Range("AF10").Value = 97.496614224
$ws.Range("AG10").Value = 116.082207546
$ws.Range("AH10").Value = 134.850822292
$ws.Range("AF11").Value = 0
$ws.Range("AG11").Value = 0
$ws.Range("AH11").Value = 0
$ws.Range("AF12").Value = 0
$ws.Range("AG12").Value = 0
$ws.Range("AH12").Value = 0
$ws.Range("AF13").Value = 55.8808590596
$ws.Range("AG13").Value = 69.97903018229999
$ws.Range("AH13").Value = 87.1682587289
$ws.Range("AF14").Value = 0
$ws.Range("AG14").Value = 0
$ws.Range("AH14").Value = 0
$ws.Range("AF15").Value = 0
$ws.Range("AG15").Value = 0
$ws.Range("AH15").Value = 0
$ws.Range("AF16").Value = 0
$ws.Range("AG16").Value = 0
$ws.Range("AH16").Value = 0
$ws.Range("AF17").Value = 71.51811290969999
$ws.Range("AG17").Value = 83.4898792445
$ws.Range("AH17").Value = 94.83043631629999
$ws.Range("AF18").Value = 74.51979790576
$ws.Range("AG18").Value = 86.63974060347
$ws.Range("AH18").Value = 0
$ws.Range("AF19").Value = 0
$ws.Range("AG19").Value = 0
$ws.Range("AH19").Value = 0
$ws.Range("AF20").Value = 0
$ws.Range("AG20").Value = 0
$ws.Range("AH20").Value = 0
$ws.Range("AF21").Value = 0
$ws.Range("AG21").Value = 0
$ws.Range("AH21").Value = 0
$ws.Range("AF22").Value = 60.851160999
$ws.Range("AG22").Value = 75.42470776499999
$ws.Range("AH22").Value = 85.755232096
$ws.Range("AF23").Value = 113.053870372
$ws.Range("AG23").Value = 136.9123724667
$ws.Range("AH23").Value = 156.0485786097
$ws.Range("AF24").Value = 75.196236946
$ws.Range("AG24").Value = 89.062112793
$ws.Range("AH24").Value = 106.651700849
$ws.Range("AF25").Value = 0
$ws.Range("AG25").Value = 0
$ws.Range("AH25").Value = 0
$ws.Range("AF26").Value = 0
$ws.Range("AG26").Value = 0
$ws.Range("AH26").Value = 0
$ws.Range("AF27").Value = 0
$ws.Range("AG27").Value = 0
$ws.Range("AH27").Value = 0
$ws.Range("AF28").Value = 75.513032749
$ws.Range("AG28").Value = 94.486585324
$ws.Range("AH28").Value = 115.689266338
$ws.Range("AF29").Value = 43.7300336234
$ws.Range("AG29").Value = 68.1845941214
$ws.Range("AH29").Value = 85.9006444444
$ws.Range("AF30").Value = 57.8923168888
$ws.Range("AG30").Value = 79.8039617264
$ws.Range("AH30").Value = 97.59974628339999
$ws.Range("AF31").Value = 80.48688252620001
$ws.Range("AG31").Value = 81.5736840917
$ws.Range("AH31").Value = 83.6186979527
$ws.Range("AF32").Value = 81.55335022200001
$ws.Range("AG32").Value = 93.657725305
$ws.Range("AH32").Value = 109.372537696
$ws.Range("AF33").Value = 0
$ws.Range("AG33").Value = 0
$ws.Range("AH33").Value = 0
$ws.Range("AF34").Value = 106.71652432
$ws.Range("AG34").Value = 129.559830121
$ws.Range("AH34").Value = 155.295083402
$ws.Range("AF35").Value = 0
$ws.Range("AG35").Value = 0
$ws.Range("AH35").Value = 0
$ws.Range("AF36").Value = 117.346722634
$ws.Range("AG36").Value = 134.538752376
$ws.Range("AH36").Value = 0
$ws.Range("AF37").Value = 0
$ws.Range("AG37").Value = 0
$ws.Range("AH37").Value = 0
$ws.Range("AF38").Value = 82.65475698100001
$ws.Range("AG38").Value = 107.271065678
$ws.Range("AH38").Value = 130.592298752
$ws.Range("AF39").Value = 91.5006471005
$ws.Range("AG39").Value = 109.0869068761
$ws.Range("AH39").Value = 127.5684101063
$ws.Range("AF40").Value = 55.375485867
$ws.Range("AG40").Value = 55.375485867
$ws.Range("AH40").Value = 55.375485867
$ws.Range("AF41").Value = 0
$ws.Range("AG41").Value = 0
$ws.Range("AH41").Value = 0
$ws.Range("AF42").Value = 79.102341307
$ws.Range("AG42").Value = 105.3188281953
$ws.Range("AH42").Value = 0
$ws.Range("AF43").Value = 0
$ws.Range("AG43").Value = 0
$ws.Range("AH43").Value = 0
$ws.Range("AF44").Value = 71.281038159
$ws.Range("AG44").Value = 88.412876766
$ws.Range("AH44").Value = 103.163733472
$ws.Range("AF45").Value = 52.705811749
$ws.Range("AG45").Value = 63.016966618
$ws.Range("AH45").Value = 73.402168129
$ws.Range("AF46").Value = 51.872209149
$ws.Range("AG46").Value = 66.26112639500001
$ws.Range("AH46").Value = 80.42873317900001
$ws.Range("AF47").Value = 75.556626194
$ws.Range("AG47").Value = 94.40515643099999
$ws.Range("AH47").Value = 109.044053422
$ws.Range("AF48").Value = 97.624565797
$ws.Range("AG48").Value = 113.303162893
$ws.Range("AH48").Value = 137.560560446
$ws.Range("AF49").Value = 0
$ws.Range("AG49").Value = 0
$ws.Range("AH49").Value = 0
$ws.Range("AF50").Value = 72.61480394500001
$ws.Range("AG50").Value = 0
$ws.Range("AH50").Value = 0
$ws.Range("AF51").Value = 60.252062814
$ws.Range("AG51").Value = 76.023184211
$ws.Range("AH51").Value = 0
$ws.Range("AF52").Value = 83.365981737
$ws.Range("AG52").Value = 101.725452821
$ws.Range("AH52").Value = 0
$ws.Range("AF53").Value = 71.552019311
$ws.Range("AG53").Value = 92.298918266
$ws.Range("AH53").Value = 0
$ws.Range("AF54").Value = 93.87450332909999
$ws.Range("AG54").Value = 115.232602715
$ws.Range("AH54").Value = 128.1463967038
$ws.Range("AF55").Value = 75.85890205699999
$ws.Range("AG55").Value = 0
$ws.Range("AH55").Value = 0
$ws.Range("AF56").Value = 118.687446339
$ws.Range("AG56").Value = 135.011023018
$ws.Range("AH56").Value = 151.616284415
$ws.Range("AF57").Value = 0
$ws.Range("AG57").Value = 0
$ws.Range("AH57").Value = 0
$ws.Range("AF58").Value = 52.447994867
$ws.Range("AG58").Value = 0
$ws.Range("AH58").Value = 0
$ws.Range("AF59").Value = 85.908890035
$ws.Range("AG59").Value = 104.010235133
$ws.Range("AH59").Value = 0
$ws.Range("AF60").Value = 0
$ws.Range("AG60").Value = 0
$ws.Range("AH60").Value = 0
$ws.Range("AF61").Value = 0
$ws.Range("AG61").Value = 0
$ws.Range("AH61").Value = 0
$ws.Range("AF62").Value = 85.88790569299999
$ws.Range("AG62").Value = 95.810562464
$ws.Range("AH62").Value = 99.15457831800001
$ws.Range("AF63").Value = 74.8552253232
$ws.Range("AG63").Value = 91.9242211235
$ws.Range("AH63").Value = 107.9400711688
$ws.Range("AF64").Value = 0
$ws.Range("AG64").Value = 0
$ws.Range("AH64").Value = 0
$ws.Range("AF65").Value = 58.623882499
$ws.Range("AG65").Value = 76.571656111
$ws.Range("AH65").Value = 93.370780098
$ws.Range("AF66").Value = 0
$ws.Range("AG66").Value = 0
$ws.Range("AH66").Value = 0
$ws.Range("AF67").Value = 64.94725393
$ws.Range("AG67").Value = 80.330597814
$ws.Range("AH67").Value = 94.596163937
$ws.Range("AF68").Value = 0
$ws.Range("AG68").Value = 0
$ws.Range("AH68").Value = 0
$ws.Range("AF69").Value = 127.480112912
$ws.Range("AG69").Value = 162.091526286
$ws.Range("AH69").Value = 194.844537196
$ws.Range("AF70").Value = 70.9797245953
$ws.Range("AG70").Value = 0
$ws.Range("AH70").Value = 0
$ws.Range("AF71").Value = 95.82696878260001
$ws.Range("AG71").Value = 119.295739353
$ws.Range("AH71").Value = 141.0301472726
$ws.Range("AF72").Value = 0
$ws.Range("AG72").Value = 0
$ws.Range("AH72").Value = 0
$ws.Range("AF73").Value = 0
$ws.Range("AG73").Value = 0
$ws.Range("AH73").Value = 0
$ws.Range("AF74").Value = 101.916427861
$ws.Range("AG74").Value = 107.917648646
$ws.Range("AH74").Value = 127.79899625
$ws.Range("AF75").Value = 101.903279745
$ws.Range("AG75").Value = 123.107824233
$ws.Range("AH75").Value = 137.905860639
$ws.Range("AF76").Value = 0
$ws.Range("AG76").Value = 0
$ws.Range("AH76").Value = 0
$ws.Range("AF77").Value = 0
$ws.Range("AG77").Value = 0
$ws.Range("AH77").Value = 0
$ws.Range("AF78").Value = 106.924648905
$ws.Range("AG78").Value = 138.208306918
$ws.Range("AH78").Value = 158.930296279
